$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ------------------------------------------------------------------
# 1) Reorder underlying industry data: reverse columns B:M so each
#    column holds the series that used to be in the mirrored slot.
#    Column A (Year) is untouched.
# ------------------------------------------------------------------
# Header row
$ws.Cells.Item(1, 2).Value = "Unclassified"
$ws.Cells.Item(1, 3).Value = "Trade, transportation, and utilities"
$ws.Cells.Item(1, 4).Value = "Public administration"
$ws.Cells.Item(1, 5).Value = "Professional and business services"
$ws.Cells.Item(1, 6).Value = "Other services"
$ws.Cells.Item(1, 7).Value = "Natural resources and mining"
$ws.Cells.Item(1, 8).Value = "Manufacturing"
$ws.Cells.Item(1, 9).Value = "Leisure and hospitality"
$ws.Cells.Item(1, 10).Value = "Information"
$ws.Cells.Item(1, 11).Value = "Financial activities"
$ws.Cells.Item(1, 12).Value = "Education and health services"
$ws.Cells.Item(1, 13).Value = "Construction"

# Data row 2 (year 2018)
$ws.Cells.Item(2, 2).Value = 0.006805035730522748
$ws.Cells.Item(2, 3).Value = 19.54292844543958
$ws.Cells.Item(2, 4).Value = 2.365884125300287
$ws.Cells.Item(2, 5).Value = 22.7821254531684
$ws.Cells.Item(2, 6).Value = 3.255075424433371
$ws.Cells.Item(2, 7).Value = 0.3595326828550536
$ws.Cells.Item(2, 8).Value = 6.935465540657757
$ws.Cells.Item(2, 9).Value = 15.62776455514549
$ws.Cells.Item(2, 10).Value = 0.7916524491822863
$ws.Cells.Item(2, 11).Value = 7.764545759674239
$ws.Cells.Item(2, 12).Value = 16.55892027760536
$ws.Cells.Item(2, 13).Value = 4.009300250807671

# Data row 3 (year 2019)
$ws.Cells.Item(3, 2).Value = 0.009982364488121698
$ws.Cells.Item(3, 3).Value = 19.10069987221597
$ws.Cells.Item(3, 4).Value = 2.328109228951937
$ws.Cells.Item(3, 5).Value = 22.84297740365181
$ws.Cells.Item(3, 6).Value = 3.279761309904424
$ws.Cells.Item(3, 7).Value = 0.3793298080747728
$ws.Cells.Item(3, 8).Value = 7.11631673545411
$ws.Cells.Item(3, 9).Value = 15.84423074586871
$ws.Cells.Item(3, 10).Value = 0.8074624036498024
$ws.Cells.Item(3, 11).Value = 7.364766705807896
$ws.Cells.Item(3, 12).Value = 16.95338235566001
$ws.Cells.Item(3, 13).Value = 3.972981066272426

# Data row 4 (year 2020)
$ws.Cells.Item(4, 2).Value = 0.01894498474031892
$ws.Cells.Item(4, 3).Value = 19.38428749575421
$ws.Cells.Item(4, 4).Value = 2.580072203940768
$ws.Cells.Item(4, 5).Value = 23.09750755480284
$ws.Cells.Item(4, 6).Value = 2.956604151933394
$ws.Cells.Item(4, 7).Value = 0.4191581134878389
$ws.Cells.Item(4, 8).Value = 7.208572613810943
$ws.Cells.Item(4, 9).Value = 14.21703866358321
$ws.Cells.Item(4, 10).Value = 0.7471434753702665
$ws.Cells.Item(4, 11).Value = 7.49629975769393
$ws.Cells.Item(4, 12).Value = 17.81185247074892
$ws.Cells.Item(4, 13).Value = 4.062518514133349

# Data row 5 (year 2021)
$ws.Cells.Item(5, 2).Value = 0.03629480923652448
$ws.Cells.Item(5, 3).Value = 19.43470914341789
$ws.Cells.Item(5, 4).Value = 2.08014335579231
$ws.Cells.Item(5, 5).Value = 21.63167626397817
$ws.Cells.Item(5, 6).Value = 2.971633379384585
$ws.Cells.Item(5, 7).Value = 0.4094502243818375
$ws.Cells.Item(5, 8).Value = 7.095625325377129
$ws.Cells.Item(5, 9).Value = 15.52735151289121
$ws.Cells.Item(5, 10).Value = 0.7338346534481274
$ws.Cells.Item(5, 11).Value = 7.940612655677051
$ws.Cells.Item(5, 12).Value = 18.16098995065114
$ws.Cells.Item(5, 13).Value = 3.977678725764017

# Data row 6 (year 2022)
$ws.Cells.Item(6, 2).Value = 0.05946399993439822
$ws.Cells.Item(6, 3).Value = 18.35274018711996
$ws.Cells.Item(6, 4).Value = 2.392613497360425
$ws.Cells.Item(6, 5).Value = 20.98537184282759
$ws.Cells.Item(6, 6).Value = 3.040230068110822
$ws.Cells.Item(6, 7).Value = 0.434627499520511
$ws.Cells.Item(6, 8).Value = 7.186489792071732
$ws.Cells.Item(6, 9).Value = 16.21744346431808
$ws.Cells.Item(6, 10).Value = 0.8087097991078164
$ws.Cells.Item(6, 11).Value = 8.332522490807408
$ws.Cells.Item(6, 12).Value = 18.32138646308894
$ws.Cells.Item(6, 13).Value = 3.868400895732309

# ------------------------------------------------------------------
# 2) Recolor each chart series with the new custom palette and
#    refresh the cached series names/colors to match the new data.
# ------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart

$chart.SeriesCollection(1).Format.Fill.ForeColor.RGB = 16559976  # 68AFFC
$chart.SeriesCollection(2).Format.Fill.ForeColor.RGB = 10892098  # 4233A6
$chart.SeriesCollection(3).Format.Fill.ForeColor.RGB = 14542213  # 85E5DD
$chart.SeriesCollection(4).Format.Fill.ForeColor.RGB = 6711338  # 2A6866
$chart.SeriesCollection(5).Format.Fill.ForeColor.RGB = 7921254  # 66DE78
$chart.SeriesCollection(6).Format.Fill.ForeColor.RGB = 5084949  # 15974D
$chart.SeriesCollection(7).Format.Fill.ForeColor.RGB = 7393716  # B4D170
$chart.SeriesCollection(8).Format.Fill.ForeColor.RGB = 15464  # 683C00
$chart.SeriesCollection(9).Format.Fill.ForeColor.RGB = 5537482  # CA7E54
$chart.SeriesCollection(10).Format.Fill.ForeColor.RGB = 4726658  # 821F48
$chart.SeriesCollection(11).Format.Fill.ForeColor.RGB = 6839286  # F65B68
$chart.SeriesCollection(12).Format.Fill.ForeColor.RGB = 13356779  # EBCECB

# ------------------------------------------------------------------
# 3) Remove the major gridlines from the value (percent) axis.
# ------------------------------------------------------------------
$valAx = $chart.Axes(2, 1)
$valAx.HasMajorGridlines = $false

# ------------------------------------------------------------------
# 4) Move/resize the chart: from the far right of the sheet to just
#    below/left of the data table.
# ------------------------------------------------------------------
$co.Left = 0
$co.Top = 105
$co.Width = 734.9296875
$co.Height = 432

Write-Host "Done."
